# Changes of 23rd May 2022
$wb = $excel.ActiveWorkbook

$wsCreation   = $wb.Worksheets.Item("RTECreation")
$wsSearch     = $wb.Worksheets.Item("SearchRTE")
$wsRate       = $wb.Worksheets.Item("Rate")
$wsRoute      = $wb.Worksheets.Item("RouteDetail")
$wsShipment   = $wb.Worksheets.Item("ShipmentDetails")
$wsLocJob     = $wb.Worksheets.Item("LocJob")

# RTECreation
$wsCreation.Range("C2").Value = "126156775"
$wsCreation.Range("C3").Value = "126156797"

# SearchRTE
$wsSearch.Range("A2").Value = "126156775"
$wsSearch.Range("B2").Value = "32411206"
$wsSearch.Range("C2").Value = "3416444"
$wsSearch.Range("D2").Value = "126156786"
$wsSearch.Range("A3").Value = "126156797"
$wsSearch.Range("B3").Value = "32411207"
$wsSearch.Range("C3").Value = "3416445"
$wsSearch.Range("D3").Value = "126156801"

# Rate
$wsRate.Range("A2").Value = "126156797"
$wsRate.Range("B2").Value = "3416445"
$wsRate.Range("D2").Value = "$520.00"
$wsRate.Range("E2").Value = "764 Minute(s)"
$wsRate.Range("F2").Value = "Scheduler"
$wsRate.Range("H2").Value = "Total :US$520.00"
$wsRate.Range("I2").Value = "US$520.00"
$wsRate.Range("K2").Value = "US$520.00"

# RouteDetail
$wsRoute.Range("A2").Value = "126156797"
$wsRoute.Range("B2").Value = "3416445"
$wsRoute.Range("F2").Value = "Test company order 1, 3625 Willowbend Blvd., Suite 132, Chemistry Lab Room P058, HOUSTON, TX, 77056, USA"
$wsRoute.Range("G2").Value = ""
$wsRoute.Range("J2").Value = "2022-05-23 00:30"
$wsRoute.Range("G3").Value = ""
$wsRoute.Range("H3").Value = "824"
$wsRoute.Range("I3").Value = "764"
$wsRoute.Range("J3").Value = "2022-05-23 13:14"

# ShipmentDetails
$wsShipment.Range("A2").Value = "3416445"
$wsShipment.Range("B2").Value = "126156801"

# LocJob
$wsLocJob.Range("A2").Value = "3416737"
